# Casework Contact details amend
#
# 1. Split the "Contact: <Casework Officer>" paragraph into two paragraphs,
#    so "Contact: " and "<Casework Officer>" each sit on their own line.
# 2. Remove the "Direct Dial: " label, leaving just "<Casework Officer Number>".
# 3. Remove the "Email: " label, leaving just "<Casework Officer Email>".
# 4. Move the "_GoBack" bookmark from the end of the "<Completion Date>"
#    paragraph to the start of the (former) "Email:" paragraph.

$d = $word.ActiveDocument

# 1. Break the paragraph right after "Contact: " so the officer name moves
#    onto its own paragraph.
$r = $d.Content
$r.Find.Execute("Contact: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.InsertParagraphAfter()

# 2. Drop the "Direct Dial: " label entirely.
$r2 = $d.Content
$r2.Find.Execute("Direct Dial: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 3. Remember where the "Email: " label is before wiping it, so we can drop
#    the bookmark at the start of what remains of that paragraph.
$r3 = $d.Content
$r3.Find.Execute("Email: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$emailPara = $r3.Paragraphs(1)
$emailStart = $emailPara.Range.Start

# 4. Relocate the "_GoBack" bookmark to the start of that paragraph.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
$markRange = $d.Range($emailStart, $emailStart)
$d.Bookmarks.Add("_GoBack", $markRange)
